# RoundTube_Options.xlsx - "Finish off round tube options (all 154)"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold the header row (OD / Wall / ID)
$ws.Range("A1:C1").Font.Bold = $true

# New OD/Wall rows 102-155 (Outer diameter, Wall thickness); the existing
# ID column already holds "=A<r>-2*B<r>" formulas (shared) for rows 102-150
# and recalculates automatically once A/B are populated.
$data = @(
    @(102, 4, 0.125),
    @(103, 1.625, 0.125),
    @(104, 1.5, 0.375),
    @(105, 1, 0.1875),
    @(106, 1, 0.25),
    @(107, 3, 0.1875),
    @(108, 2.5, 0.25),
    @(109, 2.75, 0.25),
    @(110, 4.5, 0.125),
    @(111, 3, 0.065),
    @(112, 1.375, 0.25),
    @(113, 5, 0.125),
    @(114, 3, 0.25),
    @(115, 8, 0.125),
    @(116, 3.5, 0.1875),
    @(117, 2.75, 0.065),
    @(118, 3, 0.75),
    @(119, 3.5, 0.25),
    @(120, 2.5, 0.375),
    @(121, 1.125, 0.25),
    @(122, 6, 0.125),
    @(123, 2.5, 0.5),
    @(124, 4, 0.1875),
    @(125, 3, 0.375),
    @(126, 3.5, 0.065),
    @(127, 5, 0.1875),
    @(128, 4, 0.25),
    @(129, 7, 0.125),
    @(130, 5, 0.25),
    @(131, 4.5, 0.25),
    @(132, 4, 0.065),
    @(133, 5, 0.065),
    @(134, 6, 0.25),
    @(135, 6, 0.1875),
    @(136, 3.5, 0.5),
    @(137, 3.5, 0.375),
    @(138, 5.5, 0.25),
    @(139, 4.5, 0.375),
    @(140, 4, 0.5),
    @(141, 7, 0.25),
    @(142, 4.25, 0.5),
    @(143, 6, 0.5),
    @(144, 4.5, 0.5),
    @(145, 4.5, 0.75),
    @(146, 5, 0.5),
    @(147, 7, 0.25),
    @(148, 7, 0.5),
    @(149, 8, 0.5),
    @(150, 12, 0.25),
    @(151, 7.5, 0.5),
    @(152, 9, 0.5),
    @(153, 9, 1),
    @(154, 8, 1.5),
    @(155, 8.5, 0.5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $od = $entry[1]
    $wall = $entry[2]
    $ws.Cells.Item($r, 1).Value = $od
    $ws.Cells.Item($r, 2).Value = $wall
}

# Rows 151-155 are brand new and need the ID formula added explicitly
$ws.Range("C151:C155").Formula = "=A151-2*B151"

# Page was set to print in portrait orientation
$ws.PageSetup.Orientation = 1

# Leave selection where the author ended up after entering the last row
$ws.Range("B156").Select() | Out-Null
